# Update column G ("K") values in row 2-18 to reflect regenerated save_data
# (diff only touches column G values; no formula/formatting changes)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 5
    3  = 6
    4  = 2
    5  = 4
    6  = 3
    7  = 4
    8  = 2
    9  = 2
    10 = 6
    11 = 5
    12 = 3
    13 = 5
    14 = 5
    15 = 3
    16 = 5
    17 = 3
    18 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
